$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new description text for Carson's row (C4) - this appends a new
# shared string and references it from C4, keeping the existing style (s="1").
$ws.Range("C4").Value = "Created the shopping cart page (p4), Created the Dairy Aisle and the dairy products (p2 and p3), and created the back page edit product page (p8). Also worked on CSS relating to the pages created (created a few classes to make the receipt display properly for example)."

# Move the active selection from C14 to C4 to match the saved cursor position.
$ws.Range("C4").Select() | Out-Null
